# Update the "Metadata" sheet (StructureDefinition summary table)
$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aca-health-insurance-oversight-system-product"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet (StructureDefinition element table)
$wsElem = $wb.Worksheets.Item("Elements")

# Extension.url row: Fixed Value column picks up the same URL change
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aca-health-insurance-oversight-system-product"

# Extension row: Constraint(s) column is cleared (the constraint text moved
# down to the Extension.extension row instead)
$wsElem.Range("AI2").Value = ""
